# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows (250 and 251) into the
# "Pepino ensalada" sheet, pushing the existing rows 250-265 down to
# 252-267.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 250; Excel shifts
# everything from 250 downward (old 250-265 become 252-267).
$ws.Rows("250:251").Insert()

# --- New row 250 ---
$ws.Range("A250").Value = 5
$ws.Range("B250").Value = "Macroferia Regional de Talca"
$ws.Range("C250").Value = "Maule"
$ws.Range("D250").Value = 44516
$ws.Range("E250").Value = 7
$ws.Range("F250").Value = 100112043
$ws.Range("G250").Value = "Pepino ensalada"
$ws.Range("H250").Value = "Alaska"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 100
$ws.Range("K250").Value = 15000
$ws.Range("L250").Value = 15000
$ws.Range("M250").Value = 15000
$ws.Range("N250").Value = "$/caja 60 unidades"
$ws.Range("O250").Value = "Región del Maule"
$ws.Range("P250").Value = 250
$ws.Range("Q250").Value = 60
$ws.Range("R250").Value = "Hortaliza"

# --- New row 251 ---
$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = 44516
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100112043
$ws.Range("G251").Value = "Pepino ensalada"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 500
$ws.Range("K251").Value = 8000
$ws.Range("L251").Value = 8000
$ws.Range("M251").Value = 8000
$ws.Range("N251").Value = "$/caja 80 unidades"
$ws.Range("O251").Value = "Región del Maule"
$ws.Range("P251").Value = 100
$ws.Range("Q251").Value = 80
$ws.Range("R251").Value = "Hortaliza"

Write-Host "Inserted rows 250-251; sheet now spans to row $($ws.UsedRange.Rows.Count)"
